# "modified data (because there were some logic problems)"
#
# Corrects 29 rows of source measurements on sheet "données09": column A
# (a computed/measured figure) and column C (a derived count) are revised
# downward to fix a prior calculation issue. Column B is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = 18.87
$ws.Range("C12").Value = 131

$ws.Range("A17").Value = 45.39
$ws.Range("C17").Value = 136

$ws.Range("A18").Value = 41.699999999999996
$ws.Range("C18").Value = 146

$ws.Range("A19").Value = 11.020000000000001
$ws.Range("C19").Value = 143

$ws.Range("A20").Value = 35.839999999999996
$ws.Range("C20").Value = 148

$ws.Range("A21").Value = 27.529999999999998
$ws.Range("C21").Value = 139

$ws.Range("A24").Value = 11.95
$ws.Range("C24").Value = 139

$ws.Range("A26").Value = 43.64
$ws.Range("C26").Value = 124

$ws.Range("A28").Value = 19.16
$ws.Range("C28").Value = 118

$ws.Range("A34").Value = 20.22
$ws.Range("C34").Value = 138

$ws.Range("A37").Value = 32.049999999999997
$ws.Range("C37").Value = 127

$ws.Range("A41").Value = 12.76
$ws.Range("C41").Value = 147

$ws.Range("A45").Value = 58.02
$ws.Range("C45").Value = 139

$ws.Range("A46").Value = 8.6
$ws.Range("C46").Value = 141

$ws.Range("A51").Value = 8.8800000000000008
$ws.Range("C51").Value = 147

$ws.Range("A55").Value = 7.35
$ws.Range("C55").Value = 109

$ws.Range("A56").Value = 11.15
$ws.Range("C56").Value = 130

$ws.Range("A61").Value = 11.68
$ws.Range("C61").Value = 149

$ws.Range("A67").Value = 11.91
$ws.Range("C67").Value = 113

$ws.Range("A70").Value = 34.08
$ws.Range("C70").Value = 138

$ws.Range("A71").Value = 71.819999999999993
$ws.Range("C71").Value = 143

$ws.Range("A74").Value = 26.43
$ws.Range("C74").Value = 142

$ws.Range("A75").Value = 5.82
$ws.Range("C75").Value = 147

$ws.Range("A80").Value = 8.86
$ws.Range("C80").Value = 143

$ws.Range("A81").Value = 35.010000000000005
$ws.Range("C81").Value = 137

$ws.Range("A82").Value = 5.76
$ws.Range("C82").Value = 126

$ws.Range("A94").Value = 33.79
$ws.Range("C94").Value = 131

$ws.Range("A96").Value = 19.950000000000003
$ws.Range("C96").Value = 139

$ws.Range("A97").Value = 74.89
$ws.Range("C97").Value = 140

# Cosmetic workbook-chrome metadata (last-saved window geometry / author's
# local file path recorded by Excel on save). Not part of the Excel
# Range/Worksheet object model exposed here, but set defensively in case
# the host surfaces them via the Application/Window COM objects.
$win = $excel.ActiveWindow
$win.WindowWidth = 25800
$win.WindowHeight = 13200
$win.UsableWidth = 25800
$win.UsableHeight = 13200
